$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "Execute" -> "Run" inside the quoted button captions (8 spots)
# -----------------------------------------------------------------
$quoteOpen  = [char]8220
$quoteClose = [char]8221

$d.Content.Find.Execute(
    ($quoteOpen + "Execute step by step" + $quoteClose), $false, $false, $false,
    $false, $false, $true, 1, $false,
    ($quoteOpen + "Run step by step" + $quoteClose), 2) | Out-Null

$d.Content.Find.Execute(
    ($quoteOpen + "Execute" + $quoteClose), $false, $false, $false,
    $false, $false, $true, 1, $false,
    ($quoteOpen + "Run" + $quoteClose), 2) | Out-Null

# -----------------------------------------------------------------
# 2) Insert a space between "#4" and the en-dash in the sprint-task
#    heading: "task #4-Testing GUI..." -> "task #4 -Testing GUI..."
# -----------------------------------------------------------------
$enDash = [char]8211
$d.Content.Find.Execute(
    ("task #4" + $enDash), $false, $false, $false,
    $false, $false, $true, 1, $false,
    ("task #4 " + $enDash), 2) | Out-Null

Write-Host "Phase 1 done"

# -----------------------------------------------------------------
# 3) Append four new rows (test cases 9-12) to the second table
# -----------------------------------------------------------------
$t2 = $d.Tables.Item(2)

$t2.Rows.Add() | Out-Null
$row = $t2.Rows.Count
$t2.Cell($row, 1).Range.Text = "9"
$t2.Cell($row, 2).Range.Text = "Choose one agent in each region"
$t2.Cell($row, 3).Range.Text = "A click to choose agent position"
$t2.Cell($row, 4).Range.Text = "An agent is marked in each region"

$t2.Rows.Add() | Out-Null
$row = $t2.Rows.Count
$t2.Cell($row, 1).Range.Text = "10"
$t2.Cell($row, 2).Range.Text = "The number of agent is greater than 1 and less than the number of spaces."
$t2.Cell($row, 3).Range.Text = ("Multiple clicks to choose agent" + [char]8217 + "s position")
$t2.Cell($row, 4).Range.Text = "Each region has many agents"

$t2.Rows.Add() | Out-Null
$row = $t2.Rows.Count
$t2.Cell($row, 1).Range.Text = "11"
$t2.Cell($row, 2).Range.Text = "One region has no agent"
$t2.Cell($row, 3).Range.Text = "None"
$t2.Cell($row, 4).Range.Text = "Error. Each region must have at least one agent."

$t2.Rows.Add() | Out-Null
$row = $t2.Rows.Count
$t2.Cell($row, 1).Range.Text = "12"
$t2.Cell($row, 2).Range.Text = "The number of agent is more than the half of spaces in a region"
$t2.Cell($row, 3).Range.Text = ("Multiple clicks to choose agent" + [char]8217 + "s position")
$t2.Cell($row, 4).Range.Text = "Error. The number of agents should less than half of open spaces in each region"

Write-Host "Phase 2 done (rows added)"

# -----------------------------------------------------------------
# 4) Move the "_GoBack" bookmark into the newly typed text of row 12
#    (it lands mid-word, right after "hal" in "half", mirroring the
#    author's last live edit position)
# -----------------------------------------------------------------
$lastCellRng = $t2.Cell($row, 2).Range
$lastCellRng.Find.Execute("hal", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackPos = $lastCellRng.End

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$goBackRng = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRng) | Out-Null

Write-Host "Phase 3 done (_GoBack moved)"
